$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 419, shifting existing rows 419:435 down to 420:436
$ws.Rows.Item(419).Insert()

# Fill the new row 419 with the new weekly record (same market/product metadata,
# new date + price observations)
$ws.Cells.Item(419, 1).Value = 5
$ws.Cells.Item(419, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(419, 3).Value = "Maule"
$ws.Cells.Item(419, 4).Value = 45075
$ws.Cells.Item(419, 5).Value = 7
$ws.Cells.Item(419, 6).Value = 100112009
$ws.Cells.Item(419, 7).Value = "Acelga"
$ws.Cells.Item(419, 8).Value = "Sin especificar"
$ws.Cells.Item(419, 9).Value = "Primera"
$ws.Cells.Item(419, 10).Value = 500
$ws.Cells.Item(419, 11).Value = 1800
$ws.Cells.Item(419, 12).Value = 1800
$ws.Cells.Item(419, 13).Value = 1800
$ws.Cells.Item(419, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(419, 15).Value = "Región del Maule"
$ws.Cells.Item(419, 16).Value = 450
$ws.Cells.Item(419, 17).Value = 4
$ws.Cells.Item(419, 18).Value = "Hortaliza"
